$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SeznamModulu")

$ws.Range("B45").Value = "Bunkr"
$ws.Range("B46").Value = "Tábor"
$ws.Range("A55").Value = "Aktualizace: 16. 10. 2017"

$ws.Activate() | Out-Null
$ws.Range("A4").Select() | Out-Null
